$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 45, shifting existing rows 45-57 down to 46-58.
$ws.Rows.Item(45).Insert()

# Populate the newly inserted row 45 with the new weekly price record.
$ws.Cells.Item(45, 1).Value = 1
$ws.Cells.Item(45, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(45, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(45, 4).Value = 45016
$ws.Cells.Item(45, 5).Value = 15
$ws.Cells.Item(45, 6).Value = 100112052
$ws.Cells.Item(45, 7).Value = "Albahaca"
$ws.Cells.Item(45, 8).Value = "Sin especificar"
$ws.Cells.Item(45, 9).Value = "Primera"
$ws.Cells.Item(45, 10).Value = 300
$ws.Cells.Item(45, 11).Value = 800
$ws.Cells.Item(45, 12).Value = 1000
$ws.Cells.Item(45, 13).Value = 900
$ws.Cells.Item(45, 14).Value = "$/paquete"
$ws.Cells.Item(45, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(45, 16).Value = 900
$ws.Cells.Item(45, 17).Value = 1
$ws.Cells.Item(45, 18).Value = "Hortaliza"
